$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45189 -> 45190) for every data row (rows 2 through 199).
$ws.Range("C2:C199").Value = 45190
